$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1001
$ws.Range("F4").Value = 1232
$ws.Range("F5").Value = 52
$ws.Range("F8").Value = 4587
$ws.Range("F9").Value = 592
$ws.Range("F11").Value = 1757
$ws.Range("F13").Value = 696
$ws.Range("F14").Value = 32
$ws.Range("F16").Value = 394
$ws.Range("F17").Value = 1111
$ws.Range("F19").Value = 794
$ws.Range("F24").Value = 122
$ws.Range("F28").Value = 2486
$ws.Range("F30").Value = 1501
$ws.Range("F31").Value = 480
$ws.Range("F32").Value = 15
$ws.Range("F34").Value = 4162

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 211
$ws.Range("F7").Value = 16
$ws.Range("F11").Value = 399
$ws.Range("F12").Value = 356

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1316
$ws.Range("F5").Value = 1718
$ws.Range("F8").Value = 189

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1316
$ws.Range("F3").Value = 1718
$ws.Range("F6").Value = 1001
$ws.Range("F7").Value = 1232
$ws.Range("F9").Value = 52
$ws.Range("F11").Value = 189
$ws.Range("F14").Value = 4587
$ws.Range("F15").Value = 592
$ws.Range("F17").Value = 1757
$ws.Range("F18").Value = 696
$ws.Range("F19").Value = 356
$ws.Range("F22").Value = 394
$ws.Range("F25").Value = 794
$ws.Range("F29").Value = 122
$ws.Range("F38").Value = 2486
$ws.Range("F43").Value = 1501
$ws.Range("F44").Value = 480
$ws.Range("F45").Value = 15
$ws.Range("F48").Value = 4162
